# Apply crypto price/volume updates (rows 2-51, columns D & E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.261.00"
$ws.Range("E2").Value = "  -3.31%  "
$ws.Range("D3").Value = "2.358.79"
$ws.Range("E3").Value = "  -3.04%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'556.01"
$ws.Range("E5").Value = "  -3.27%  "
$ws.Range("D6").Value = "'136.66"
$ws.Range("E6").Value = "  -2.88%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").Value = "2.352.56"
$ws.Range("E9").Value = "  -2.72%  "
$ws.Range("E10").Value = "  -4.21%  "
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("D12").Value = "'5.05"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("D14").Value = "'25.31"
$ws.Range("E14").Value = "  -3.26%  "
$ws.Range("D15").Value = "2.781.57"
$ws.Range("E15").Value = "  -3.85%  "
$ws.Range("E16").Value = "  -4.20%  "
$ws.Range("D17").Value = "59.234.87"
$ws.Range("E17").Value = "  -3.30%  "
$ws.Range("D18").Value = "2.345.50"
$ws.Range("E18").Value = "  -3.57%  "
$ws.Range("D19").Value = "'7.98"
$ws.Range("E19").Value = "  +10.19%  "
$ws.Range("D20").Value = "'10.38"
$ws.Range("E20").Value = "  -2.03%  "
$ws.Range("D21").Value = "'319.90"
$ws.Range("D22").Value = "'4.01"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("E23").Value = "  -1.46%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  -6.63%  "
$ws.Range("D26").Value = "'63.90"
$ws.Range("D27").Value = "'549.07"
$ws.Range("E27").Value = "  -4.02%  "
$ws.Range("D28").Value = "'8.01"
$ws.Range("E28").Value = "  -9.94%  "
$ws.Range("D29").Value = "2.470.88"
$ws.Range("D30").Value = "0.0₃0905"
$ws.Range("E30").Value = "  -1.11%  "
$ws.Range("D31").Value = "'7.90"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("E32").Value = "  -4.36%  "
$ws.Range("E33").Value = "  -4.98%  "
$ws.Range("E34").Value = "  -2.77%  "
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").Value = "'150.77"
$ws.Range("E37").Value = "  -0.73%  "
$ws.Range("D38").Value = "'0.364"
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("E39").Value = "  -2.48%  "
$ws.Range("D40").Value = "'18.04"
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("D41").Value = "'4.98"
$ws.Range("E41").Value = "  -2.56%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "'41.08"
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("D44").Value = "'1.62"
$ws.Range("E44").Value = "  -2.45%  "
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("D46").Value = "0.0₆0283"
$ws.Range("E46").Value = "  -2.12%  "
$ws.Range("D47").Value = "'138.33"
$ws.Range("E47").Value = "  -1.91%  "
$ws.Range("E48").Value = "  -1.39%  "
$ws.Range("D49").Value = "'0.580"
$ws.Range("E49").Value = "  -2.10%  "
$ws.Range("D50").Value = "'0.0496"
$ws.Range("E50").Value = "  -2.16%  "
$ws.Range("D51").Value = "'18.90"
$ws.Range("E51").Value = "  -3.06%  "
